$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the data values in B2:F2
$ws.Range("B2").Value = 762957.95475691999
$ws.Range("C2").Value = 229725.89961944899
$ws.Range("D2").Value = 59192.511157275803
$ws.Range("E2").Value = 18436.239720580299
$ws.Range("F2").Value = 6450.4833894153599

# Change number format on the data cells (B2:F4) from scientific to fixed 3-decimal
$ws.Range("B2:F4").NumberFormat = "0.000"

# Column width changes: B:C narrower (10.57), D:F stay at 11.29
$ws.Range("B:C").ColumnWidth = 10.5703125
$ws.Range("D:F").ColumnWidth = 11.28515625

# Update the selection
$ws.Range("H14").Select()
